# Updated symbol list on Wed Jan 25 12:51:45 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "303.04"
Set-TextValue "E2" "-3.53%"
Set-TextValue "D3" "35.49"
Set-TextValue "E3" "0.96%"
Set-TextValue "D4" "5.041"
Set-TextValue "E4" "-1.47%"
Set-TextValue "D5" "0.08036"
Set-TextValue "E5" "-1.52%"
Set-TextValue "D6" "1.931"
Set-TextValue "E6" "-8.67%"
Set-TextValue "D7" "7.801"
Set-TextValue "E7" "-1.97%"
Set-TextValue "B8" "GateToken"
Set-TextValue "C8" "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextValue "D8" "4.048"
Set-TextValue "E8" "-2.42%"
Set-TextValue "B9" "BTSEToken"
Set-TextValue "C9" "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextValue "D9" "2.991"
Set-TextValue "E9" "7.73%"
Set-TextValue "B10" "MXToken"
Set-TextValue "C10" "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue "D10" "0.9253"
Set-TextValue "E10" "-0.25%"
Set-TextValue "B11" "LiechtensteinCryptoassetsExchange"
Set-TextValue "C11" "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue "D11" "0.1235"
Set-TextValue "E11" "19.88%"
Set-TextValue "B12" "WazirX"
Set-TextValue "C12" "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue "D12" "0.1850"
Set-TextValue "E12" "-1.98%"
Set-TextValue "B13" "MandalaExchangeToken"
Set-TextValue "C13" "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue "D13" "0.09449"
Set-TextValue "E13" "4.64%"
Set-TextValue "B14" "BitrueCoin"
Set-TextValue "C14" "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue "D14" "0.03463"
Set-TextValue "E14" "-3.96%"
Set-TextValue "B15" "BitMartToken"
Set-TextValue "C15" "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue "D15" "0.09877"
Set-TextValue "E15" "-0.25%"
Set-TextValue "B16" "BitForexToken"
Set-TextValue "C16" "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue "D16" "0.001386"
Set-TextValue "E16" "-3.02%"
Set-TextValue "B17" "TigerCash"
Set-TextValue "C17" "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue "D17" "0.005758"
Set-TextValue "E17" "1.00%"
Set-TextValue "B18" "LEO"
Set-TextValue "C18" "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue "D18" "3.507"
Set-TextValue "E18" "0.96%"
Set-TextValue "E19" "-0.64%"
Set-TextValue "D20" "0.1283"
Set-TextValue "E20" "-4.13%"
Set-TextValue "D21" "5.052"
Set-TextValue "E21" "-0.93%"
Set-TextValue "D22" "0.2397"
Set-TextValue "E22" "8.46%"
Set-TextValue "D23" "0.04486"
Set-TextValue "E23" "-0.53%"
Set-TextValue "D24" "0.001210"
Set-TextValue "E24" "-2.78%"
Set-TextValue "D25" "0.004812"
Set-TextValue "E25" "2.21%"
Set-TextValue "D26" "0.0001249"
Set-TextValue "E26" "-0.11%"
Set-TextValue "D27" "0.0002998"
Set-TextValue "E27" "-33.33%"
Set-TextValue "D39" "0.01925"
Set-TextValue "E39" "-2.24%"
Set-TextValue "D40" "0.04749"
Set-TextValue "E40" "-2.76%"
Set-TextValue "D41" "0.007360"
Set-TextValue "E41" "-3.85%"
Set-TextValue "D42" "0.009639"
Set-TextValue "E42" "23.03%"
Set-TextValue "D43" "0.1332"
Set-TextValue "E43" "-4.28%"
Set-TextValue "D44" "0.002108"
Set-TextValue "E44" "-1.97%"
Set-TextValue "D45" "0.01085"
Set-TextValue "E45" "-7.98%"
Set-TextValue "D46" "0.00006255"
Set-TextValue "E46" "-6.80%"
Set-TextValue "D47" "0.00000000750"
Set-TextValue "E47" "0.01%"
Set-TextValue "E48" "66.08%"
Set-TextValue "E49" "-12.34%"
Set-TextValue "D50" "0.00002099"
Set-TextValue "E50" "0.01%"
Set-TextValue "D51" "0.0001999"
Set-TextValue "E51" "0.01%"
